$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need NumberFormat forced to
# Text ("@") before assignment, otherwise Excel auto-converts the string into a
# numeric value (losing formatting like "1.00" and the inline-string type).
$textForceCells = @("D5", "D6", "D10", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D24", "D25", "D27", "D29", "D31", "D32", "D33", "D34", "D37", "D38", "D39", "D40", "D41", "D43", "D45", "D46", "D47", "D48", "D49")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (price + 1h volume change), row by row in sheet order.
$ws.Range("D2").Value = "65.307.70"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "3.189.03"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "616.69"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("D6").Value = "148.83"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.179.65"
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "0.154"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("D12").Value = "0.478"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").Value = "36.19"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("D15").Value = "3.711.55"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("E16").Value = "  +3.23%  "
$ws.Range("D17").Value = "65.236.20"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "3.180.20"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "6.96"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").Value = "484.87"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").Value = "14.81"
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("D22").Value = "0.726"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("E23").Value = "  +2.97%  "
$ws.Range("D24").Value = "13.97"
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").Value = "84.85"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "8.76"
$ws.Range("E27").Value = "  +2.38%  "
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").Value = "7.11"
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("E30").Value = "  -4.88%  "
$ws.Range("D31").Value = "2.13"
$ws.Range("E31").Value = "  -4.92%  "
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").Value = "2.74"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "26.87"
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("E35").Value = "  +2.96%  "
$ws.Range("E36").Value = "  +6.07%  "
$ws.Range("D37").Value = "6.09"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "3.25"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").Value = "53.30"
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("D40").Value = "470.17"
$ws.Range("E40").Value = "  +4.08%  "
$ws.Range("D41").Value = "0.0404"
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("D43").Value = "8.43"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").Value = "2.870.21"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").Value = "2.37"
$ws.Range("E45").Value = "  +2.01%  "
$ws.Range("D46").Value = "0.272"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").Value = "2.49"
$ws.Range("E47").Value = "  +7.31%  "
$ws.Range("D48").Value = "37.80"
$ws.Range("E48").Value = "  +12.51%  "
$ws.Range("D49").Value = "27.02"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("E51").Value = "  -0.28%  "
